$d = $word.ActiveDocument

# --- Helpers ---------------------------------------------------------------

function Xml-Escape($s) {
    return ($s -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;')
}

# Locate the next occurrence of $searchText at/after character position
# $fromPos within the document and return a fresh Range spanning exactly
# the matched text (not the paragraph mark).
function Find-RangeAfter($searchText, $fromPos) {
    $scan = $d.Range($fromPos, $d.Content.End)
    $found = $scan.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    return $d.Range($scan.Start, $scan.End)
}

# Replace the contents of a Range with $newText, wrapping the text in a
# single run whose run properties (e.g. bold/italic) are given by the raw
# OOXML fragment $rPrXml (or $null for none). Using InsertXML (instead of
# Range.Text= or Find.Execute replace) keeps sibling runs - including
# zero-length <w:r/> placeholder runs - untouched.
function Replace-RangeXml($range, $rPrXml, $newText) {
    $escaped = Xml-Escape $newText
    $preserve = ""
    if ($newText -ne $newText.Trim()) { $preserve = ' xml:space="preserve"' }
    $runInner = ""
    if ($rPrXml) { $runInner += $rPrXml }
    $runInner += "<w:t$preserve>$escaped</w:t>"
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Replace-Text($searchText, $newText, $rPrXml, $fromPos) {
    $r = Find-RangeAfter $searchText $fromPos
    Replace-RangeXml $r $rPrXml $newText
    return $r.Start + $newText.Length
}

# --- Edits -------------------------------------------------------------

$pos = 0

# Heading1 title (plain run, no rPr)
$pos = Replace-Text "Play Free Book of Kings Slot - Big Wins Await" "Play Book of Kings for Free - Slot Game Review" $null $pos

# "What we like" bullet list (plain runs, no rPr)
$pos2 = Replace-Text "Free bonuses and bonus features" "Offers free bonuses and special symbols" $null $pos
$pos2 = Replace-Text "Wild symbols to replace all game symbols" "Visually pleasing graphics with an Ancient Egypt theme" $null $pos2
$pos2 = Replace-Text "Highly rewarding protagonist explorer symbol" "Mobile-compatible with no loss in quality" $null $pos2
$pos2 = Replace-Text "Engaging mobile gameplay with no loss in graphics quality" "High RTP and volatility for high prizes" $null $pos2

# "What we don't like" bullet list (plain runs, no rPr)
$pos2 = Replace-Text "High volatility may not appeal to all players" "Game symbols may not offer as high payouts" $null $pos2
$pos2 = Replace-Text "Limited paylines may not provide enough variety" "Lower frequency of wins compared to other games" $null $pos2

# Bold run near the bottom (second occurrence of the title text)
$pos2 = Replace-Text "Play Free Book of Kings Slot - Big Wins Await" "Play Book of Kings for Free - Slot Game Review" "<w:rPr><w:b/></w:rPr>" $pos2

# Italic summary paragraph at the end
$pos2 = Replace-Text "Discover the mysteries of Ancient Egypt and play Book of Kings, an online slot game with free bonuses, high rewards, and mobile compatibility. Play free now." "Discover the gameplay, graphics, symbols, and mobile compatibility of Book of Kings in our review. Play for free!" "<w:rPr><w:i/></w:rPr>" $pos2

Write-Host "Done"
